$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F18").Value = "use restrictions"
$ws.Range("F31").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F37").Value = "32_physical_and_chemical_hazards"
$ws.Range("F57").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F81").Value = "use restrictions"
$ws.Range("F92").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F95").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F101").Value = "application instructions"
$ws.Range("F111").Value = "use restrictions"
$ws.Range("F124").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F127").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F131").Value = "32_physical_and_chemical_hazards"
$ws.Range("F153").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F157").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F160").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F164").Value = "32_physical_and_chemical_hazards"
$ws.Range("F166").Value = "application instructions"
$ws.Range("F175").Value = "use restrictions"
$ws.Range("F279").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F283").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F290").Value = "application instructions || env warning - species"
$ws.Range("F291").Value = "env warning - water"
$ws.Range("F292").Value = "use restrictions || env warning - water"
$ws.Range("F298").Value = "application instructions"
$ws.Range("F299").Value = "application instructions"
$ws.Range("F300").Value = "application instructions"
$ws.Range("F303").Value = "154_pesticide_storage"
$ws.Range("F310").Value = "use restrictions"
$ws.Range("F311").Value = "use restrictions || env warning - water"
$ws.Range("F312").Value = "use restrictions"
$ws.Range("F313").Value = "use restrictions"
$ws.Range("F315").Value = "use restrictions || off target movement"
$ws.Range("F317").Value = "use restrictions || off target movement"
$ws.Range("F318").Value = "use restrictions || off target movement"
$ws.Range("F320").Value = "use restrictions || off target movement"
$ws.Range("F321").Value = "use restrictions || off target movement"
$ws.Range("F323").Value = "application instructions"
$ws.Range("F324").Value = "application instructions"
$ws.Range("F325").Value = "application instructions"
$ws.Range("F326").Value = "application instructions"
$ws.Range("F327").Value = "application instructions"
$ws.Range("F328").Value = "application instructions"
$ws.Range("F330").Value = "application instructions"
$ws.Range("F332").Value = "application instructions"
$ws.Range("F333").Value = "application instructions"
$ws.Range("F336").Value = "use restrictions"
$ws.Range("F341").Value = "application instructions"
$ws.Range("F367").Value = "safety procedures"
$ws.Range("F368").Value = "safety procedures"
$ws.Range("F370").Value = "mixing"
$ws.Range("F371").Value = "mixing"
$ws.Range("F372").Value = "mixing"
$ws.Range("F373").Value = "mixing"
$ws.Range("F375").Value = "application instructions"
$ws.Range("F377").Value = "mixing"
$ws.Range("F379").Value = "mixing"
$ws.Range("F381").Value = "mixing"
$ws.Range("F383").Value = "mixing"
$ws.Range("F385").Value = "mixing"
$ws.Range("F391").Value = "mixing"
$ws.Range("F393").Value = "mixing"
$ws.Range("F394").Value = "mixing"
$ws.Range("F395").Value = "mixing"
$ws.Range("F396").Value = "mixing"
$ws.Range("F398").Value = "mixing"
$ws.Range("F401").Value = "application instructions || off target movement"
$ws.Range("F403").Value = "application instructions"
$ws.Range("F415").Value = "application instructions"
$ws.Range("F416").Value = "application instructions"
$ws.Range("F419").Value = "application instructions"
$ws.Range("F420").Value = "application instructions"
$ws.Range("F421").Value = "application instructions"
$ws.Range("F424").Value = "use restrictions"
$ws.Range("F425").Value = "use restrictions"
$ws.Range("F428").Value = "application instructions"
$ws.Range("F430").Value = "application instructions"
$ws.Range("F440").Value = "application instructions"
$ws.Range("F451").Value = "application instructions"
$ws.Range("F453").Value = "application instructions"
$ws.Range("F456").Value = "application instructions"
$ws.Range("F465").Value = "application instructions"
$ws.Range("F471").Value = "use restrictions || application instructions || mixing"
$ws.Range("F472").Value = "application instructions"
$ws.Range("F473").Value = "mixing"
$ws.Range("F474").Value = "mixing"
$ws.Range("F475").Value = "mixing"
$ws.Range("F476").Value = "mixing"
$ws.Range("F479").Value = "mixing"
$ws.Range("F480").Value = "off target movement"
$ws.Range("F483").Value = "application instructions"
$ws.Range("F484").Value = "use restrictions || application instructions"
$ws.Range("F486").Value = "application instructions"
$ws.Range("F487").Value = "use restrictions || application instructions"
$ws.Range("F491").Value = "use restrictions || application instructions || mixing"
$ws.Range("F492").Value = "use restrictions"
$ws.Range("F493").Value = "use restrictions || application instructions"
$ws.Range("F494").Value = "mixing"
$ws.Range("F496").Value = "mixing"
$ws.Range("F499").Value = "mixing"
$ws.Range("F500").Value = "mixing"
$ws.Range("F501").Value = "mixing"
$ws.Range("F503").Value = "mixing || application instructions"
$ws.Range("F504").Value = "mixing"
$ws.Range("F505").Value = "mixing"
$ws.Range("F506").Value = "mixing"
$ws.Range("F507").Value = "mixing"
$ws.Range("F508").Value = "mixing"
$ws.Range("F511").Value = "mixing || application instructions"
$ws.Range("F513").Value = "mixing"
$ws.Range("F515").Value = "mixing"
$ws.Range("F526").Value = "mixing"
$ws.Range("F528").Value = "use restrictions"
$ws.Range("F531").Value = "use restrictions"
$ws.Range("F532").Value = "use restrictions"
